$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 23 values (B23: 25 -> 37, C23: 2 -> 5)
$ws.Range("B23").Value = 37
$ws.Range("C23").Value = 5

# Add new row 24 data
$ws.Range("B24").Value = 9
$ws.Range("C24").Value = 4

# Fill down the formulas from row 23 into the newly-populated row 24
$ws.Range("D23").AutoFill($ws.Range("D23:D24"))
$ws.Range("F23:L23").AutoFill($ws.Range("F23:L24"))
$ws.Range("N23:T23").AutoFill($ws.Range("N23:T24"))

# Move the selection back to A1 on the data sheet
$ws.Range("A1").Select()

# Restore the window width recorded the last time the workbook was saved
$excel.ActiveWindow.Width = 23895
